$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newline = [char]10

$ws.Range("B20").Value = "Group I: <15" + $newline + "Group II: <7"
$ws.Range("D20").Value = "Group I: Som <150,000" + $newline + "Group II: Som <230,000"
$ws.Range("B21").Value = "Group I: 15-50" + $newline + "Group II: 7-50"
$ws.Range("D21").Value = "Group I: Som 150,000 - 500,000" + $newline + "Group II: Som 230,000 - 500,000"
$ws.Range("B22").Value = "Group I: 51-200" + $newline + "Group II: 16-50"
$ws.Range("B23").Value = "Group I: >200" + $newline + "Group II: >50"
